$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CERN")

# Row 4 - Inventory
$ws.Range("B4").Value = 23000000.0
$ws.Range("C4").Value = 16000000.0
$ws.Range("D4").Value = 19000000.0
$ws.Range("E4").Value = 23000000.0
$ws.Range("F4").Value = 23000000.0

# Row 13 - Accounts Payable
$ws.Range("B13").Value = 236000000.0
$ws.Range("C13").Value = 256000000.0
$ws.Range("D13").Value = 262000000.0
$ws.Range("E13").Value = 296000000.0
$ws.Range("F13").Value = 273000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = 376000000.0
$ws.Range("C21").Value = 392000000.0
$ws.Range("D21").Value = 382000000.0
$ws.Range("E21").Value = 379000000.0
$ws.Range("F21").Value = 378000000.0
